# Applies the "Added Fresh-Rotten dataset- 10-07-19" edit:
#  - Rename the original sheet "Feuil1" -> "Fruit Recognition"
#  - Append two rows (52, 53) with a single value (3) in column A
#  - Remove the hidden "_xlchart.v1.*" defined names
#  - Add a new worksheet "Fresh-Rotten" with a small header + 2 data rows
#  - Leave "Fresh-Rotten" as the active/selected sheet

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename + extra rows ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Fruit Recognition"

$ws1.Range("A52").Value = 3
$ws1.Range("A53").Value = 3

# Matches the saved selection/view state recorded in the target file.
[void]$ws1.Range("B52:I53").Select()

# --- Workbook-level: drop the hidden chart helper names ---------------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- Sheet 2: new "Fresh-Rotten" worksheet, appended after sheet 1 ----------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Fresh-Rotten"

$headers = @("Number of labels", "Number of layers", "Learning rate", "Batch Size", "Epochs", "Accuracy", "Train Time (s)", "Used Transform")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws2.Range("A2").Value = 6
$ws2.Range("B2").Value = 6
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 50
$ws2.Range("E2").Value = 3
$ws2.Range("F2").Value = 0.86174944403261677
$ws2.Range("G2").Value = 1993
$ws2.Range("H2").Value = $false

$ws2.Range("A3").Value = 6
$ws2.Range("B3").Value = 6
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 50
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 0.8528539659006672
$ws2.Range("G3").Value = 403
$ws2.Range("H3").Value = $false

# Select/activate to match the recorded view state, and leave this as the active tab.
[void]$ws2.Range("A2:H3").Select()
[void]$ws2.Activate()
